# Applies the "Bug Jagd" checklist update:
#  - Adds a new comment string "TODO: Einmal organisiert durchgehen"
#  - Tags rows 7, 8 and 10 (column F) with that comment
#  - Flips the "IST" value of row 10 (Variablenhandhabung) from 1 to 0
#  - Scrolls the sheet view back to the top and moves the selection to G5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$comment = "TODO: Einmal organisiert durchgehen"

$ws.Range("F7").Value = $comment
$ws.Range("F8").Value = $comment
$ws.Range("F10").Value = $comment

$ws.Range("D10").Value = 0

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("G5").Select()
